$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45177) for every data row
# (rows 2-410). Update it to 45178 (one day later) for all of them.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 410 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45178
